# Update the "summary" sheet so the automation code can read test case
# settings (test_id, description, browser_type, test_type, result) for
# each test case, including a new test003 entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("summary")

# Header row
$ws.Range("A1").Value = "test_id"
$ws.Range("B1").Value = "description"
$ws.Range("C1").Value = "browser_type"
$ws.Range("D1").Value = "test_type"
$ws.Range("E1").Value = "result"

# Existing test cases, now carrying browser_type/test_type settings
$ws.Range("A2").Value = "test001"
$ws.Range("B2").Value = "Test register"
$ws.Range("C2").Value = "firefox"
$ws.Range("D2").Value = "auto"

$ws.Range("A3").Value = "test002"
$ws.Range("B3").Value = "Test login"
$ws.Range("C3").Value = "firefox"
$ws.Range("D3").Value = "auto"

$ws.Range("A4").Value = "google"
$ws.Range("B4").Value = "Test search text google"
$ws.Range("C4").Value = "firefox"
$ws.Range("D4").Value = "auto"

# New test case row
$ws.Range("A5").Value = "test003"
$ws.Range("B5").Value = "Test something"
$ws.Range("C5").Value = "firefox"
$ws.Range("D5").Value = "manual"

# Make "summary" the active/selected sheet (was "google"), with K3 selected
$ws.Activate()
$ws.Range("K3").Select()
